# Fix the ASSISTS (F) column: cells were stored as text ("inline string")
# but should be numeric. Re-assigning a numeric value converts the cell
# type to numeric while keeping the same value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$assists = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 2
    22 = 2
    23 = 2
    24 = 2
    25 = 2
    26 = 2
    27 = 2
    28 = 2
    29 = 2
    30 = 2
    31 = 2
    32 = 2
    33 = 2
    34 = 2
    35 = 2
    36 = 3
    37 = 3
    38 = 3
    39 = 3
    40 = 3
    41 = 3
}

foreach ($row in $assists.Keys) {
    $ws.Cells.Item($row, 6).Value = $assists[$row]
}

# Fix the CHAMPION (H) column: several rows incorrectly listed other
# champions; they should all read "Yone".
$champRows = @(5, 11, 17, 23, 25, 29, 35, 37, 41)
foreach ($row in $champRows) {
    $ws.Cells.Item($row, 8).Value = "Yone"
}
